$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Metadata sheet: bump the "Date" row value (B8) to the new commit date.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# ---------------------------------------------------------------------------
# 2) Elements sheet: the two trailing "Mapping" columns swap places - column
#    AK held "Mapping: RIM Mapping" and AL held "Mapping: Spécification
#    métier vers l'extension ROR ContactTelecomConfidentialityLevel"; after
#    the edit the business-spec mapping comes first (AK) and RIM Mapping
#    comes second (AL). Each row's data travels with its column, so we just
#    swap the AK/AL cell contents row by row (only rows whose two values
#    actually differ need touching).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Row 1 - headers
$ws.Cells.Item(1, 37).Value = "Mapping: Spécification métier vers l'extension ROR ContactTelecomConfidentialityLevel"
$ws.Cells.Item(1, 38).Value = "Mapping: RIM Mapping"

# Row 3 - was AK3="n/a", AL3="" -> now AK3="", AL3="n/a"
$ws.Cells.Item(3, 37).Value = ""
$ws.Cells.Item(3, 38).Value = "n/a"

# Row 5 - was AK5="N/A", AL5="" -> now AK5="", AL5="N/A"
$ws.Cells.Item(5, 37).Value = ""
$ws.Cells.Item(5, 38).Value = "N/A"

# Row 6 - was AK6="N/A", AL6="niveauConfidentialite" -> now AK6="niveauConfidentialite", AL6="N/A"
$ws.Cells.Item(6, 37).Value = "niveauConfidentialite"
$ws.Cells.Item(6, 38).Value = "N/A"

# Column widths follow the data: the (now wider) business-spec mapping text
# sits in AK, the (now narrower) RIM mapping column sits in AL.
$ws.Columns.Item(37).ColumnWidth = 91.5
$ws.Columns.Item(38).ColumnWidth = 24.15
